$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.350.92'

$ws.Range("D3").Value = '2.324.46'
$ws.Range("E3").Value = '  +0.79%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("E7").Value = '  -0.85%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  -0.62%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.66'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.70%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.72'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.77%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0798'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.00%  '

$ws.Range("E13").Value = '  +0.37%  '

$ws.Range("E14").Value = '  +1.75%  '

$ws.Range("D15").Value = '2.689.36'
$ws.Range("E15").Value = '  +0.96%  '

$ws.Range("D16").Value = '2.323.31'
$ws.Range("E16").Value = '  +0.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.793'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.97%  '

$ws.Range("D18").Value = '43.253.61'
$ws.Range("E18").Value = '  +0.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.93%  '

$ws.Range("D20").Value = '0.0₃0901'
$ws.Range("E20").Value = '  -0.49%  '

$ws.Range("E21").Value = '  +0.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.47'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.36%  '

$ws.Range("E24").Value = '  +4.35%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.46'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.05%  '

$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.61%  '

$ws.Range("E28").Value = '  +0.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '164.36'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.43%  '

$ws.Range("E30").Value = '  +0.71%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.22'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.13%  '

$ws.Range("E33").Value = '  -0.36%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.92'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0706'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.00%  '

$ws.Range("E37").Value = '  -1.33%  '

$ws.Range("E38").Value = '  -0.18%  '

$ws.Range("E39").Value = '  +0.97%  '

$ws.Range("E40").Value = '  +1.93%  '

$ws.Range("E41").Value = '  -0.46%  '

$ws.Range("D42").Value = '1.986.85'
$ws.Range("E42").Value = '  -1.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.14'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.93%  '

$ws.Range("E44").Value = '  +5.53%  '

$ws.Range("E45").Value = '  -0.40%  '

$ws.Range("E46").Value = '  -0.81%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.81'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.51%  '

$ws.Range("D48").Value = '2.555.89'
$ws.Range("E48").Value = '  +0.98%  '

$ws.Range("B49").Value = 'HuobiToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.87'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.91%  '

$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.98'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.49%  '

$ws.Range("E51").Value = '  +0.64%  '
